$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'69.713.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.77%  "

# Row 3
$ws.Range("D3").Value = "'3.503.99"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.09%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "'604.27"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.99%  "

# Row 6
$ws.Range("D6").Value = "'170.96"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.88%  "

# Row 7
$ws.Range("D7").Value = "'0.615"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.39%  "

# Row 8
$ws.Range("D8").Value = "'3.498.63"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.05%  "

# Row 9
$ws.Range("E9").Value = "  +0.00%  "

# Row 10
$ws.Range("D10").Value = "'0.200"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +5.51%  "

# Row 11
$ws.Range("D11").Value = "'6.70"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.67%  "

# Row 12
$ws.Range("E12").Value = "  -2.75%  "

# Row 13
$ws.Range("E13").Value = "  -0.27%  "

# Row 14
$ws.Range("D14").Value = "'0.0000279"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.96%  "

# Row 15
$ws.Range("D15").Value = "'4.076.62"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.16%  "

# Row 16
$ws.Range("D16").Value = "'620.57"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -8.42%  "

# Row 17
$ws.Range("E17").Value = "  -4.04%  "

# Row 18
$ws.Range("D18").Value = "'3.503.67"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.04%  "

# Row 19
$ws.Range("D19").Value = "'69.708.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.81%  "

# Row 20
$ws.Range("E20").Value = "  -2.19%  "

# Row 21
$ws.Range("D21").Value = "'17.27"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.24%  "

# Row 22
$ws.Range("D22").Value = "'0.882"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.54%  "

# Row 23
$ws.Range("D23").Value = "'9.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -11.69%  "

# Row 24
$ws.Range("D24").Value = "'15.78"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.38%  "

# Row 25
$ws.Range("D25").Value = "'96.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.00%  "

# Row 26
$ws.Range("E26").Value = "  -0.79%  "

# Row 27
$ws.Range("E27").Value = "  -0.09%  "

# Row 28
$ws.Range("E28").Value = "  -2.74%  "

# Row 29
$ws.Range("E29").Value = "  -3.25%  "

# Row 30
$ws.Range("D30").Value = "'33.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.43%  "

# Row 31
$ws.Range("D31").Value = "'8.38"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.08%  "

# Row 32
$ws.Range("E32").Value = "  -4.96%  "

# Row 33
$ws.Range("E33").Value = "  -1.40%  "

# Row 34
$ws.Range("D34").Value = "'6.94"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.62%  "

# Row 35
$ws.Range("D35").Value = "'561.95"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.47%  "

# Row 36
$ws.Range("D36").Value = "'10.73"
$ws.Range("D36").Style = "Normal"

# Row 37
$ws.Range("D37").Value = "'3.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.70%  "

# Row 38
$ws.Range("D38").Value = "'56.99"
$ws.Range("D38").Style = "Normal"

# Row 39
$ws.Range("E39").Value = "  -4.04%  "

# Row 40
$ws.Range("D40").Value = "'0.999"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.06%  "

# Row 41
$ws.Range("D41").Value = "'0.142"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.93%  "

# Row 42
$ws.Range("D42").Value = "'0.0447"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.40%  "

# Row 43
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "'0.325"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.73%  "

# Row 44
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "'3.322.59"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.95%  "

# Row 45
$ws.Range("D45").Value = "'0.0₃0707"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.37%  "

# Row 46
$ws.Range("E46").Value = "  +1.85%  "

# Row 47
$ws.Range("D47").Value = "'32.89"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.69%  "

# Row 48
$ws.Range("D48").Value = "'2.61"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.24%  "

# Row 49
$ws.Range("E49").Value = "  -3.36%  "

# Row 50
$ws.Range("D50").Value = "'134.95"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.17%  "

# Row 51
$ws.Range("E51").Value = "  -1.38%  "
